$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.107348918914795
$ws.Range("B1").Value = 2.239099979400635
$ws.Range("C1").Value = 10.05295944213867
$ws.Range("D1").Value = 1.334017992019653
$ws.Range("E1").Value = 1.271888732910156
